$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (RandomForestClassifier)
$ws.Range("C2").Value = 0.7593989071038252
$ws.Range("D2").Value = 0.707074990093247
$ws.Range("E2").Value = 0.8117813694304328
$ws.Range("F2").Value = 0.8466666666666667
$ws.Range("G2").Value = 0.6721311475409836
$ws.Range("H2").Value = 0.8639455782312925
$ws.Range("I2").Value = 0.8466666666666667
$ws.Range("J2").Value = 0.8552188552188552

# Row 3 (XGBClassifier)
$ws.Range("C3").Value = 0.7140983606557377
$ws.Range("D3").Value = 0.6635254331778327
$ws.Range("E3").Value = 0.7656883460930433
$ws.Range("F3").Value = 0.92
$ws.Range("G3").Value = 0.5081967213114754
$ws.Range("H3").Value = 0.8214285714285714
$ws.Range("I3").Value = 0.92
$ws.Range("J3").Value = 0.8679245283018868

# Row 4 (LogisticRegression)
$ws.Range("C4").Value = 0.7739890710382513
$ws.Range("D4").Value = 0.7231135310556935
$ws.Range("E4").Value = 0.8249246581455041
$ws.Range("F4").Value = 0.8266666666666667
$ws.Range("G4").Value = 0.7213114754098361
$ws.Range("H4").Value = 0.8794326241134752
$ws.Range("I4").Value = 0.8266666666666667
$ws.Range("J4").Value = 0.852233676975945
